$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was recorded for "Sandia" at Femacal de La Calera.
# It belongs chronologically right after the current row 675, so insert a
# fresh row at 676; everything that used to be row 676.. shifts down to
# 677.. (dimension grows from A1:R734 to A1:R735 automatically).
$ws.Rows.Item(676).Insert()

$ws.Range("A676").Value = 3
$ws.Range("B676").Value = "Femacal de La Calera"
$ws.Range("C676").Value = "Coquimbo"
$ws.Range("D676").Value = 45223
$ws.Range("E676").Value = 5
$ws.Range("F676").Value = 100112028
$ws.Range("G676").Value = "Sandia"
$ws.Range("H676").Value = "Sin especificar"
$ws.Range("I676").Value = "Primera"
$ws.Range("J676").Value = 180
$ws.Range("K676").Value = 700
$ws.Range("L676").Value = 700
$ws.Range("M676").Value = 700
$ws.Range("N676").Value = "$/kilo (volumen en unidades)"
$ws.Range("O676").Value = "Per" + [char]0xFA
$ws.Range("P676").Value = 700
$ws.Range("Q676").Value = 1
$ws.Range("R676").Value = "Hortaliza"
